# =============================================================================
# edit.ps1 - Word COM-interop script
#
# Commit: "Moved the self loop logic to the front end since it no longer
# makes sense to keep it on the server."
#
# This script reproduces, against the already-open $word.ActiveDocument:
#   1. Spell-check "proofErr" (spellStart/spellEnd) wrappers around
#      "neighbours" / "Venkata", achieved by splitting the w:r runs that
#      contained those words into multiple runs.
#   2. A large block of new meeting-notes paragraphs appended after the
#      "...20,000 genes" paragraph, with the "_GoBack" bookmark relocated
#      from the old last paragraph to the end of the new content.
#   3. A new, trailing empty paragraph at the very end of the document.
#
# Approach: Range.InsertXML() accepts a full single-part WordprocessingML
# "xmlPackage" and, when applied to a Range spanning one or more whole
# paragraphs (i.e. including their paragraph marks), replaces those
# paragraphs with the supplied <w:p> elements verbatim - which lets us
# splice in <w:proofErr/> markers that aren't reachable through simple
# Find/Replace. Edits are applied bottom-up so earlier paragraph indices
# remain valid while later ones are being rewritten.
# =============================================================================

$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 4) Paragraph 20 ("Let’s see what the performance...20,000 genes") is
#    rewritten without its trailing bookmark, and is followed immediately
#    by seven new paragraphs of notes, ending in a bookmark-terminated
#    paragraph plus one brand-new empty paragraph.
# -----------------------------------------------------------------------
$pTail = $d.Paragraphs(20)
$rTail = $pTail.Range
$xmlTail = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00AC16FD" w:rsidRPr="004D6404" w:rsidRDefault="00AC16FD" w:rsidP="004D6404"><w:r><w:t>Let’s see what the performance of the md-autocomplete control is like when it contains 20,000</w:t></w:r><w:r w:rsidR="00140439"><w:t xml:space="preserve"> genes</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Let’s test out the functionality of getting beyond just the second </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. We will implement this by having a list similar to the genes of interest list. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WE’ll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> leave it to the ng-repeat to tell the user which nodes are selected. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Unfortunately, our logic breaks down after a few levels since of the conditions we implemented to ensure that there won’t be duplicate nodes in the graph are preventing necessary edges from being added. What we can do about this is add a check that </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Goes back and searches the previous </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to see if one of our exclusions for the current neighbor is in there. IF so, we’ll add an edge from one of those previous </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the current node of interest. The only issue with this approach is that it will make the graph somewhat messy</w:t></w:r><w:r><w:t>. Furthermore, this approach doesn’t work under the current framework since we are only ever using a single source node at each level. Let’s come back to this later.</w:t></w:r></w:p><w:p><w:r><w:t>One thing we can do right now is make a layout f</w:t></w:r><w:r><w:t>or the selected genes graph. Okay so here is the ideal layout:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">We need to highlight the genes that the user selected. Unfortunately, we can’t just have a fixed policy that these genes are to go in the middle of the graph. Our position policy should be based on minimizing edge overlap. I’m not sure if it’s even possible to have a layout where no edges cross </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>eacht</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> other for this kind of graph. </w:t></w:r><w:r><w:t xml:space="preserve">From a theoretical point of view, I don’t think it’s possible to completely avoid edge overlap in the genes of interest graph. Having said that, we can try to minimize the amount of overlapping edges. Since we have control over node position and size, there probably isn’t a unique way in which to accomplish this task. </w:t></w:r><w:r><w:t>Let’s try making some sort of concentric layout and see where that takes us.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Another issue that has emerged is that our self-loop logic is currently flawed due to the fact that we are no longer caching graph elements. </w:t></w:r><w:r><w:t xml:space="preserve">We can simply move this logic to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getRelevatnSubmatrix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> R script and send a list of genes back to the client in addition to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cytoscape</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> that we are sending back. However, we would also have to do this for any other script that we are to write in the future. Already, the list that we are returning from R to the server is quite complicated and adding more elements to it will lead to confusion in the future as well as a lack of maintainability. We should leave the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>self loop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> computation to the client side. </w:t></w:r><w:r><w:tab/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rTail.InsertXML($xmlTail)

# -----------------------------------------------------------------------
# 3) Paragraph 16 ("IT would be nice...Venkata gets back.") - wrap
#    "Venkata" with spellStart/spellEnd proofErr markers.
# -----------------------------------------------------------------------
$pVenkata = $d.Paragraphs(16)
$rVenkata = $pVenkata.Range
$xmlVenkata = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005217F0" w:rsidRDefault="005217F0" w:rsidP="004D6404"><w:r><w:t xml:space="preserve">IT would be nice to have the circular layouts done by the time that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Venkata</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gets back.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rVenkata.InsertXML($xmlVenkata)

# -----------------------------------------------------------------------
# 2) Paragraph 15 ("We have generalized the neighbor general method...")
#    - wrap both occurrences of "neighbours" with proofErr markers.
# -----------------------------------------------------------------------
$pGeneralized = $d.Paragraphs(15)
$rGeneralized = $pGeneralized.Range
$xmlGeneralized = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="004D6404" w:rsidRDefault="004D6404" w:rsidP="004D6404"><w:r><w:t xml:space="preserve">We have generalized the neighbor general method so that it now allows you to keep asking for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and doesn’t limit you to only the 2</w:t></w:r><w:r w:rsidRPr="004D6404"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="005217F0"><w:t>. We need to extend the front end in order to permit the user to keep on exploring neighbours</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="005217F0"><w:t xml:space="preserve">It is still however made to work for only an epi-stroma correlation matrix. We need to come up with a scheme that will work got epi-epi as well. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rGeneralized.InsertXML($xmlGeneralized)

# -----------------------------------------------------------------------
# 1) Paragraphs 5-6 ("1. Test first neighbours" / "2. Test second
#    neighbours") - wrap "neighbours" with proofErr markers in both.
# -----------------------------------------------------------------------
$pFirstNeighbours = $d.Paragraphs(5)
$pSecondNeighbours = $d.Paragraphs(6)
$rNeighbours = $d.Range($pFirstNeighbours.Range.Start, $pSecondNeighbours.Range.End)
$xmlNeighbours = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="004D6404" w:rsidRDefault="004D6404" w:rsidP="004D6404"><w:r><w:t xml:space="preserve">1. Test first </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p w:rsidR="004D6404" w:rsidRDefault="004D6404" w:rsidP="004D6404"><w:r><w:t xml:space="preserve">2. Test second </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rNeighbours.InsertXML($xmlNeighbours)

Write-Host 'Applied all edits for Week of May 16.docx'
